# Applies the cryptos-list price/volume refresh described in the commit message
# ("Updated cryptos list ... with GitHub Actions").
#
# Every cell in this sheet is stored as text (t="inlineStr" in the original
# workbook), including price cells that look like numbers (e.g. "1.00",
# "0.999"). A plain $range.Value = "1.00" would be auto-coerced to the number
# 1, and prefixing with a bare apostrophe keeps it text but also stamps the
# cell with a new quote-prefix style, changing its style index from the
# original 0. Re-applying the "Normal" style right after clears that stamp,
# so the cell ends up back at style 0 with the text value preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.176.69'
$ws.Range("E2").Value = '  +1.72%  '
# Row 3
$ws.Range("D3").Value = '2.382.71'
$ws.Range("E3").Value = '  +3.85%  '
# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
# Row 5
$ws.Range("D5").Value = '''303.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.95%  '
# Row 6
$ws.Range("D6").Value = '''96.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.47%  '
# Row 7
$ws.Range("E7").Value = '  +0.68%  '
# Row 8
$ws.Range("E8").Value = '  -0.15%  '
# Row 9
$ws.Range("D9").Value = '''0.501'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.21%  '
# Row 10
$ws.Range("D10").Value = '''34.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.26%  '
# Row 11
$ws.Range("D11").Value = '''0.0789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.33%  '
# Row 12
$ws.Range("E12").Value = '  +2.32%  '
# Row 13
$ws.Range("D13").Value = '''18.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.49%  '
# Row 14
$ws.Range("D14").Value = '''6.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.68%  '
# Row 15
$ws.Range("D15").Value = '2.752.34'
$ws.Range("E15").Value = '  +3.83%  '
# Row 16
$ws.Range("D16").Value = '2.367.84'
$ws.Range("E16").Value = '  +1.58%  '
# Row 17
$ws.Range("E17").Value = '  +3.99%  '
# Row 18
$ws.Range("D18").Value = '43.157.14'
$ws.Range("E18").Value = '  +1.79%  '
# Row 19
$ws.Range("E19").Value = '  +0.26%  '
# Row 20
$ws.Range("D20").Value = '''6.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.00%  '
# Row 21
$ws.Range("D21").Value = '0.0₃0889'
$ws.Range("E21").Value = '  +0.30%  '
# Row 22
$ws.Range("D22").Value = '''68.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.62%  '
# Row 23
$ws.Range("E23").Value = '  +0.75%  '
# Row 24
$ws.Range("D24").Value = '''235.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '
# Row 25
$ws.Range("D25").Value = '''2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.24%  '
# Row 26
$ws.Range("E26").Value = '  +0.00%  '
# Row 27
$ws.Range("D27").Value = '''24.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.11%  '
# Row 28
$ws.Range("E28").Value = '  +15.25%  '
# Row 29
$ws.Range("E29").Value = '  +1.48%  '
# Row 30
$ws.Range("D30").Value = '''31.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.10%  '
# Row 31
$ws.Range("D31").Value = '''0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '
# Row 32
$ws.Range("D32").Value = '''5.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.69%  '
# Row 33
$ws.Range("E33").Value = '  +6.69%  '
# Row 34
$ws.Range("D34").Value = '''17.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.41%  '
# Row 35
$ws.Range("E35").Value = '  +7.14%  '
# Row 36
$ws.Range("E36").Value = '  +3.43%  '
# Row 37
$ws.Range("E37").Value = '  -1.37%  '
# Row 38
$ws.Range("E38").Value = '  -0.99%  '
# Row 39
$ws.Range("D39").Value = '''2.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.31%  '
# Row 40
$ws.Range("D40").Value = '''22.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.91%  '
# Row 41
$ws.Range("E41").Value = '  +0.49%  '
# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '''104.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -36.65%  '
# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.956.51'
$ws.Range("E43").Value = '  +0.49%  '
# Row 44
$ws.Range("E44").Value = '  +1.19%  '
# Row 45
$ws.Range("E45").Value = '  +2.28%  '
# Row 46
$ws.Range("E46").Value = '  +1.04%  '
# Row 47
$ws.Range("D47").Value = '''9.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.70%  '
# Row 48
$ws.Range("D48").Value = '2.606.10'
$ws.Range("E48").Value = '  +3.40%  '
# Row 49
$ws.Range("D49").Value = '''52.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '
# Row 50
$ws.Range("E50").Value = '  +3.22%  '
# Row 51
$ws.Range("D51").Value = '''71.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.12%  '
